# Pais worksheet refresh: updated COVID-19 country data pull (13 Sep 2020, 19:51)
# and re-ranking swaps: Italia/Irak, Suecia/China/Marruecos, Suazilandia/Mozambique
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Septiembre de 2020 a las 19:51"

# Row 4
$ws.Cells.Item(4, 2).Value = 6694486
$ws.Cells.Item(4, 3).Value = 17885
$ws.Cells.Item(4, 4).Value = 3959252
$ws.Cells.Item(4, 5).Value = 2536945
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 161
$ws.Cells.Item(4, 8).Value = 198289

# Row 5
$ws.Cells.Item(5, 2).Value = 4837952
$ws.Cells.Item(5, 3).Value = 86164
$ws.Cells.Item(5, 4).Value = 3767542
$ws.Cells.Item(5, 5).Value = 990720
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1076
$ws.Cells.Item(5, 8).Value = 79690

# Row 6
$ws.Cells.Item(6, 2).Value = 4319184
$ws.Cells.Item(6, 3).Value = 3326
$ws.Cells.Item(6, 4).Value = 3553421
$ws.Cells.Item(6, 5).Value = 634355
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 134
$ws.Cells.Item(6, 8).Value = 131408

# Row 14
$ws.Cells.Item(14, 2).Value = 434748
$ws.Cells.Item(14, 3).Value = 2082
$ws.Cells.Item(14, 4).Value = 406326
$ws.Cells.Item(14, 5).Value = 16473
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 54
$ws.Cells.Item(14, 8).Value = 11949

# Row 16
$ws.Cells.Item(16, 2).Value = 381094
$ws.Cells.Item(16, 3).Value = 7183
$ws.Cells.Item(16, 4).Value = 89059
$ws.Cells.Item(16, 5).Value = 261119
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 6
$ws.Cells.Item(16, 8).Value = 30916

# Row 21
$ws.Cells.Item(21, 2).Value = 291162
$ws.Cells.Item(21, 3).Value = 1527
$ws.Cells.Item(21, 4).Value = 258833
$ws.Cells.Item(21, 5).Value = 25273
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 57
$ws.Cells.Item(21, 8).Value = 7056

# Row 22
$ws.Cells.Item(22, 1).Value = "Irak"
$ws.Cells.Item(22, 2).Value = 290309
$ws.Cells.Item(22, 3).Value = 3531
$ws.Cells.Item(22, 4).Value = 224705
$ws.Cells.Item(22, 5).Value = 57590
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 73
$ws.Cells.Item(22, 8).Value = 8014

# Row 23
$ws.Cells.Item(23, 1).Value = "Italia"
$ws.Cells.Item(23, 2).Value = 287753
$ws.Cells.Item(23, 3).Value = 1456
$ws.Cells.Item(23, 4).Value = 213634
$ws.Cells.Item(23, 5).Value = 38509
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 7
$ws.Cells.Item(23, 8).Value = 35610

# Row 41
$ws.Cells.Item(41, 1).Value = "Marruecos"
$ws.Cells.Item(41, 2).Value = 86686
$ws.Cells.Item(41, 3).Value = 2251
$ws.Cells.Item(41, 4).Value = 67528
$ws.Cells.Item(41, 5).Value = 17580
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 25
$ws.Cells.Item(41, 8).Value = 1578

# Row 42
$ws.Cells.Item(42, 1).Value = "Suecia"
$ws.Cells.Item(42, 2).Value = 86505
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 5846

# Row 43
$ws.Cells.Item(43, 1).Value = "China"
$ws.Cells.Item(43, 2).Value = 85184
$ws.Cells.Item(43, 3).Value = 10
$ws.Cells.Item(43, 4).Value = 80399
$ws.Cells.Item(43, 5).Value = 151
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 4634

# Row 59
$ws.Cells.Item(59, 2).Value = 48254
$ws.Cells.Item(59, 3).Value = 247
$ws.Cells.Item(59, 4).Value = 34037
$ws.Cells.Item(59, 5).Value = 12605
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 7
$ws.Cells.Item(59, 8).Value = 1612

# Row 72
$ws.Cells.Item(72, 2).Value = 30985
$ws.Cells.Item(72, 3).Value = 255
$ws.Cells.Item(72, 4).Value = 23364
$ws.Cells.Item(72, 5).Value = 5837
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 1
$ws.Cells.Item(72, 8).Value = 1784

# Row 77
$ws.Cells.Item(77, 2).Value = 24310
$ws.Cells.Item(77, 3).Value = 641
$ws.Cells.Item(77, 4).Value = 8334
$ws.Cells.Item(77, 5).Value = 15735
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 241

# Row 113
$ws.Cells.Item(113, 1).Value = "Mozambique"
$ws.Cells.Item(113, 2).Value = 5269
$ws.Cells.Item(113, 3).Value = 229
$ws.Cells.Item(113, 4).Value = 2960
$ws.Cells.Item(113, 5).Value = 2274
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 35

# Row 114
$ws.Cells.Item(114, 1).Value = "Suazilandia"
$ws.Cells.Item(114, 2).Value = 5050
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 4188
$ws.Cells.Item(114, 5).Value = 764
$ws.Cells.Item(114, 6).Value = 0
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 98

# Row 134
$ws.Cells.Item(134, 2).Value = 3234
$ws.Cells.Item(134, 3).Value = 39
$ws.Cells.Item(134, 4).Value = 2996
$ws.Cells.Item(134, 5).Value = 226
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 7).Value = 0
$ws.Cells.Item(134, 8).Value = 12

# Row 162
$ws.Cells.Item(162, 2).Value = 1319
$ws.Cells.Item(162, 3).Value = 3
$ws.Cells.Item(162, 4).Value = 1210
$ws.Cells.Item(162, 5).Value = 27
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 82
